# "Seg the new transmittal tests and marked for test execution"
#
# The "Documents_New" test case (row 2) had a stale "PASS" result recorded
# in the Results column. Re-marking it for execution means clearing that
# recorded result so the test shows as not-yet-run again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Results column for the Documents_New test case -> clear the old PASS verdict.
$ws.Range("E2").Value = ""
